$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (42 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 315.77777
$ws.Range("J17").Value = 315.77777
$ws.Range("L17").Value = 947.33331
$ws.Range("N17").Value = -1283.33331
$ws.Range("H32").Value = 5270.8184
$ws.Range("I32").Value = 5031.8335
$ws.Range("J32").Value = 5557.6
$ws.Range("K32").Value = 5031.8335
$ws.Range("L32").Value = 5557.6
$ws.Range("M32").Value = -4705.8335
$ws.Range("N32").Value = -6209.6
$ws.Range("H80").Value = 1720.9166
$ws.Range("I80").Value = 837.5
$ws.Range("K80").Value = 2512.5
$ws.Range("M80").Value = -1514.5
$ws.Range("H83").Value = 1720.9166
$ws.Range("I83").Value = 837.5
$ws.Range("K83").Value = 7537.5
$ws.Range("M83").Value = -2545.5
$ws.Range("H106").Value = 850
$ws.Range("I106").Value = 850
$ws.Range("K106").Value = 850
$ws.Range("M106").Value = -219
$ws.Range("H112").Value = 2624.7334
$ws.Range("J112").Value = 2814.3076
$ws.Range("L112").Value = 8442.9228
$ws.Range("N112").Value = -10658.9228
$ws.Range("H113").Value = 14912.8
$ws.Range("I113").Value = 16141
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 16141
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -12887
$ws.Range("N113").Value = -16508
$ws.Range("H132").Value = 1003
$ws.Range("I132").Value = 1003
$ws.Range("K132").Value = 3009
$ws.Range("M132").Value = -479
$ws.Range("H138").Value = 8579.388999999999
$ws.Range("J138").Value = 8956.154
$ws.Range("L138").Value = 26868.462
$ws.Range("N138").Value = -37148.462

# ---- Sheet: ARM (37 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4254.778
$ws.Range("I32").Value = 4005.3635
$ws.Range("K32").Value = 4005.3635
$ws.Range("M32").Value = -3718.3635
$ws.Range("H45").Value = 4999.875
$ws.Range("J45").Value = 4999.875
$ws.Range("L45").Value = 4999.875
$ws.Range("N45").Value = -5753.875
$ws.Range("H61").Value = 1072.7142
$ws.Range("I61").Value = 1099
$ws.Range("J61").Value = 1007
$ws.Range("K61").Value = 1099
$ws.Range("L61").Value = 1007
$ws.Range("M61").Value = -887
$ws.Range("N61").Value = -1431
$ws.Range("H102").Value = 2403
$ws.Range("I102").Value = 1955
$ws.Range("J102").Value = 3299
$ws.Range("K102").Value = 1955
$ws.Range("L102").Value = 3299
$ws.Range("M102").Value = -333
$ws.Range("N102").Value = -6543
$ws.Range("H122").Value = 2495
$ws.Range("I122").Value = 2490
$ws.Range("K122").Value = 7470
$ws.Range("M122").Value = -5020
$ws.Range("H132").Value = 2748
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 1072.7142
$ws.Range("I136").Value = 1099
$ws.Range("J136").Value = 1007
$ws.Range("K136").Value = 3297
$ws.Range("L136").Value = 3021
$ws.Range("M136").Value = -747
$ws.Range("N136").Value = -8121

# ---- Sheet: CRP (37 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2867.8333
$ws.Range("I31").Value = 2158.3333
$ws.Range("J31").Value = 4996.3335
$ws.Range("K31").Value = 2158.3333
$ws.Range("L31").Value = 4996.3335
$ws.Range("M31").Value = -1863.3333
$ws.Range("N31").Value = -5586.3335
$ws.Range("H34").Value = 2867.8333
$ws.Range("I34").Value = 2158.3333
$ws.Range("J34").Value = 4996.3335
$ws.Range("K34").Value = 2158.3333
$ws.Range("L34").Value = 4996.3335
$ws.Range("M34").Value = -1956.3333
$ws.Range("N34").Value = -5400.3335
$ws.Range("H87").Value = 10497
$ws.Range("I87").Value = 10497
$ws.Range("K87").Value = 10497
$ws.Range("M87").Value = -9311
$ws.Range("H90").Value = 10497
$ws.Range("I90").Value = 10497
$ws.Range("K90").Value = 31491
$ws.Range("M90").Value = -25563
$ws.Range("H105").Value = 7917.6665
$ws.Range("I105").Value = 6502.3335
$ws.Range("J105").Value = 9333
$ws.Range("K105").Value = 6502.3335
$ws.Range("L105").Value = 9333
$ws.Range("M105").Value = -4755.3335
$ws.Range("N105").Value = -12827
$ws.Range("H132").Value = 1148.8572
$ws.Range("I132").Value = 1148.8572
$ws.Range("K132").Value = 3446.5716
$ws.Range("M132").Value = -916.5715999999998
$ws.Range("H134").Value = 990
$ws.Range("I134").Value = 990
$ws.Range("K134").Value = 2970
$ws.Range("M134").Value = -435

# ---- Sheet: CUL (48 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 666666800
$ws.Range("I4").Value = 666666800
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2000000400
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2000000288
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 127.2
$ws.Range("I7").Value = 204
$ws.Range("K7").Value = 612
$ws.Range("M7").Value = -500
$ws.Range("H9").Value = 433.33334
$ws.Range("H68").Value = 4215
$ws.Range("I68").Value = 4443
$ws.Range("J68").Value = 3987
$ws.Range("K68").Value = 13329
$ws.Range("L68").Value = 11961
$ws.Range("M68").Value = -12518
$ws.Range("N68").Value = -13583
$ws.Range("H71").Value = 4215
$ws.Range("I71").Value = 4443
$ws.Range("J71").Value = 3987
$ws.Range("K71").Value = 39987
$ws.Range("L71").Value = 35883
$ws.Range("M71").Value = -35931
$ws.Range("N71").Value = -43995
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H131").Value = 977.1429000000001
$ws.Range("J131").Value = 990
$ws.Range("L131").Value = 2970
$ws.Range("N131").Value = -13050
$ws.Range("H137").Value = 2967.3845
$ws.Range("I137").Value = 2404.75
$ws.Range("J137").Value = 3217.4443
$ws.Range("K137").Value = 7214.25
$ws.Range("L137").Value = 9652.332900000001
$ws.Range("M137").Value = -2114.25
$ws.Range("N137").Value = -19852.3329
$ws.Range("H140").Value = 910.7778
$ws.Range("I140").Value = 910.7778
$ws.Range("K140").Value = 2732.3334
$ws.Range("M140").Value = 2447.6666

# ---- Sheet: GSM (15 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1286
$ws.Range("I113").Value = 1286
$ws.Range("K113").Value = 1286
$ws.Range("M113").Value = 884
$ws.Range("H122").Value = 4566.467
$ws.Range("I122").Value = 4408.727
$ws.Range("K122").Value = 13226.181
$ws.Range("M122").Value = -10776.181
$ws.Range("H132").Value = 4961.2856
$ws.Range("I132").Value = 3866.3
$ws.Range("J132").Value = 7698.75
$ws.Range("K132").Value = 11598.9
$ws.Range("L132").Value = 23096.25
$ws.Range("M132").Value = -9068.900000000001
$ws.Range("N132").Value = -28156.25

# ---- Sheet: LTW (27 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2331.5557
$ws.Range("I40").Value = 2331.5557
$ws.Range("K40").Value = 2331.5557
$ws.Range("M40").Value = -2195.5557
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1928.3334
$ws.Range("I122").Value = 1922.8
$ws.Range("K122").Value = 5768.4
$ws.Range("M122").Value = -3318.4
$ws.Range("H132").Value = 8369
$ws.Range("I132").Value = 6998
$ws.Range("J132").Value = 11111
$ws.Range("K132").Value = 20994
$ws.Range("L132").Value = 33333
$ws.Range("M132").Value = -18464
$ws.Range("N132").Value = -38393
$ws.Range("H136").Value = 7103.1113
$ws.Range("I136").Value = 7053.1875
$ws.Range("K136").Value = 21159.5625
$ws.Range("M136").Value = -18609.5625

# ---- Sheet: WVR (16 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 26
$ws.Range("I32").Value = 26
$ws.Range("K32").Value = 26
$ws.Range("M32").Value = 291
$ws.Range("H113").Value = 557
$ws.Range("I113").Value = 600.0909
$ws.Range("K113").Value = 1800.2727
$ws.Range("M113").Value = 369.7273
$ws.Range("H122").Value = 1663
$ws.Range("I122").Value = 1695.875
$ws.Range("K122").Value = 5087.625
$ws.Range("M122").Value = -2637.625
$ws.Range("H136").Value = 1110.2858
$ws.Range("I136").Value = 1104.091
$ws.Range("K136").Value = 3312.273
$ws.Range("M136").Value = -762.2729999999997
